$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spells")

# Copy the formatting of the last existing data row (98) down onto the
# new rows (99-105) so the new rows inherit the same cell styles.
$ws.Range("A98:Q98").Copy()
$ws.Range("A99:Q105").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New spell data (lore of the wild spells)
$newRows = @(
    @{ Row=99;  Name="Bray Scream";     NewRevised="New"; Level=3; School="Evocation";    Artificer="No";  Bard="Yes"; Cleric="No";  Druid="Yes"; Paladin="No"; Ranger="No";  Sorcerer="No";  Warlock="No";  Wizard="No" },
    @{ Row=100; Name="Devolve";         NewRevised="New"; Level=1; School="Enchantment";  Artificer="No";  Bard="Yes"; Cleric="No";  Druid="Yes"; Paladin="No"; Ranger="No";  Sorcerer="Yes"; Warlock="No";  Wizard="No" },
    @{ Row=101; Name="Mad Rampage";     NewRevised="New"; Level=4; School="Enchantment";  Artificer="No";  Bard="Yes"; Cleric="No";  Druid="No";  Paladin="No"; Ranger="No";  Sorcerer="Yes"; Warlock="Yes"; Wizard="No" },
    @{ Row=102; Name="Mystic Signal";   NewRevised="New"; Level=5; School="Trasmutation"; Artificer="No";  Bard="Yes"; Cleric="Yes"; Druid="Yes"; Paladin="No"; Ranger="Yes"; Sorcerer="No";  Warlock="No";  Wizard="Yes" },
    @{ Row=103; Name="Savage Dominion"; NewRevised="New"; Level=7; School="Conjuration";  Artificer="No";  Bard="No";  Cleric="No";  Druid="Yes"; Paladin="No"; Ranger="No";  Sorcerer="No";  Warlock="Yes"; Wizard="No" },
    @{ Row=104; Name="Traitor-Kin";     NewRevised="New"; Level=3; School="Enchantment";  Artificer="No";  Bard="Yes"; Cleric="No";  Druid="Yes"; Paladin="No"; Ranger="Yes"; Sorcerer="No";  Warlock="No";  Wizard="No" },
    @{ Row=105; Name="Vile Tide";       NewRevised="New"; Level=2; School="Evocation";    Artificer="No";  Bard="No";  Cleric="Yes"; Druid="Yes"; Paladin="No"; Ranger="No";  Sorcerer="Yes"; Warlock="No";  Wizard="No" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Name
    $ws.Range("B$row").Value = $r.NewRevised
    $ws.Range("C$row").Value = $r.Level
    $ws.Range("D$row").Value = $r.School
    $ws.Range("E$row").Value = $r.Artificer
    $ws.Range("F$row").Value = $r.Bard
    $ws.Range("G$row").Value = $r.Cleric
    $ws.Range("H$row").Value = $r.Druid
    $ws.Range("I$row").Value = $r.Paladin
    $ws.Range("J$row").Value = $r.Ranger
    $ws.Range("K$row").Value = $r.Sorcerer
    $ws.Range("L$row").Value = $r.Warlock
    $ws.Range("M$row").Value = $r.Wizard
    $ws.Range("N$row").Value = "0.0.0"
    $ws.Range("O$row").Value = "Complete"
    $ws.Range("P$row").Value = "Publicly Released"
    $ws.Range("Q$row").Value = "Not on website"
}

# Update the Creatures sheet's remembered selection (it's no longer the
# active tab, but its last selection moves to G295).
$wsCreatures = $wb.Worksheets.Item("Creatures")
$wsCreatures.Activate()
$wsCreatures.Range("G295").Select()

# Make the Spells sheet the active/visible tab, scrolled to the new rows,
# with the selection left on E110.
$ws.Activate()
$ws.Range("A86").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E110").Select()
